# Apply odds/value updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "G3"  = 2.38
    "J3"  = 3.1
    "S3"  = 1.5
    "T3"  = 2.5
    "AC3" = 7.5
    "AH3" = 9
    "AI3" = 15
    "AT3" = 2.5

    "Q8" = 2
    "R8" = 1.85

    "I11"  = 4.33
    "J11"  = 2.2
    "L11"  = 4.5
    "U11"  = 1.57
    "V11"  = 2.25
    "AB11" = 21
    "AK11" = 51
    "AL11" = 34
    "AM11" = 34
    "AQ11" = 23
    "BA11" = 81

    "G21"  = 5
    "I21"  = 1.53
    "L21"  = 2.05
    "O21"  = 1.14
    "P21"  = 5.5
    "U21"  = 1.62
    "V21"  = 2.2
    "Y21"  = 17
    "AA21" = 41
    "AC21" = 17
    "AG21" = 151
    "AH21" = 9.5
    "AI21" = 9
    "AK21" = 12
    "AM21" = 21
    "AO21" = 26
    "AW21" = 3.75
    "AX21" = 7.5

    "M23" = 1.08
    "N23" = 8
    "Q23" = 2.3
    "R23" = 1.6
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
